$d = $word.ActiveDocument

# Replace the huge Java stack-trace text (the orange/bold run) with the
# short message, matching the diff exactly. We use wildcard Find/Replace
# because the original text is very long; the trailing "?" wildcard also
# consumes the single literal line-break character that sits right before
# the closing </w:t>, so the replacement text has no trailing newline.
$found = $d.Content.Find.Execute(
    "aqlFeatureAccess*Main.java:1472)?",
    $true,
    $false,
    $true,
    $false,
    $false,
    $true,
    1,
    $false,
    "Feature name not found in EClass EObject",
    2
)
